# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" list is reordered from descending (1911 -> 1607) to
# ascending (1607 -> 1911) and the Salario Basico / Valor Mora columns are
# refreshed with the new period's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ascending period list: 1607..1612, 1701..1712, 1801..1812, 1901..1911
$periods = @(
    "1607","1608","1609","1610","1611","1612",
    "1701","1702","1703","1704","1705","1706","1707","1708","1709","1710","1711","1712",
    "1801","1802","1803","1804","1805","1806","1807","1808","1809","1810","1811","1812",
    "1901","1902","1903","1904","1905","1906","1907","1908","1909","1910","1911"
)

$firstRow = 16
$lastRow = 56
$newValorMora = 828116

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $firstRow + $i
    $period = $periods[$i]

    # First 26 periods (1607..1808) keep the lower "Salario Basico"; the
    # remaining periods (1809..1911) use the higher figure.
    if ($i -lt 26) {
        $salario = 27580
    } else {
        $salario = 31249
    }

    $ws.Cells.Item($row, 5).Value = $period
    $ws.Cells.Item($row, 6).Value = $salario
    $ws.Cells.Item($row, 7).Value = $newValorMora
}
